# Pay Period Report: remove the "Excused w/ Pay" and "Excused w/o Pay"
# columns (G:H). The "Totals" column (I) shifts left to become the new
# column G, along with its formatting/merges.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G:H").Delete() | Out-Null

# Restore the cursor/selection to where the author left it.
$ws.Range("L13").Select() | Out-Null
